$wb = $excel.ActiveWorkbook

# The localization status "Ready for handoff" moved on to "In Translation".
# That status string shows up on the Overview sheet (one column per locale:
# zh-cn in E2, de-de in F2) and on each locale's own detail sheet (its
# "Status" column, C2).
$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
# Shorter text -> narrower auto-fitted columns for zh-cn/de-de.
$wsOverview.Columns("E:F").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns("C:C").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
